# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column G ("K")
$kValues = @{
    2  = 9
    3  = 12
    4  = 8
    5  = 2
    6  = 8
    7  = 5
    8  = 10
    9  = 11
    10 = 8
    11 = 11
    12 = 15
    13 = 4
    14 = 8
    15 = 3
    16 = 7
    17 = 7
    18 = 5
    19 = 3
    20 = 4
    21 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
